## Add a centered "Questionnaire 28" header (Arial, 12pt) to the document,
## matching what Word does when a user enables a page header and types text.

$d = $word.ActiveDocument

$wdHeaderFooterPrimary = 1
$wdAlignParagraphCenter = 1

$section = $d.Sections(1)
$header = $section.Headers($wdHeaderFooterPrimary)

# Make sure this section's header is not just inheriting the (nonexistent)
# previous section's content, then type the questionnaire number into it.
$header.LinkToPrevious = $false
$header.Range.InsertAfter("Questionnaire 28")

# Apply the "Header" paragraph style and center it.
$paragraph = $header.Range.Paragraphs(1)
$paragraph.Style = "Header"
$paragraph.Alignment = $wdAlignParagraphCenter

# Format the run text (exclude the trailing paragraph mark so the run
# formatting doesn't leak into the paragraph mark's properties).
$textRange = $header.Range
[void]$textRange.MoveEnd(1, -1)
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
